$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Daily-effort entries added for Day 5 (column I) on the first two backlog items.
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1

# Estimated effort correction for the second backlog item.
$ws.Range("D7").Value = 4

# Leave the same selection state the author ended up with.
$ws.Range("I9").Select() | Out-Null
